$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column values (revised GDP figures)
$ws.Range("B2").Value = 63.89110288734933
$ws.Range("B3").Value = 64.65787901671419
$ws.Range("B4").Value = 65.32103675021895
$ws.Range("B5").Value = 65.59044457945527
$ws.Range("B6").Value = 66.18106943585795
$ws.Range("B7").Value = 66.95820740480886
$ws.Range("B8").Value = 68.34669390933442
$ws.Range("B9").Value = 68.67827277608683
$ws.Range("B10").Value = 70.12893031812848
$ws.Range("B11").Value = 71.37235106844989
$ws.Range("B12").Value = 73.34110058979218
$ws.Range("B14").Value = 74.52262576899044
$ws.Range("B15").Value = 74.77247468872079
$ws.Range("B16").Value = 75.76352864378279
$ws.Range("B18").Value = 76.09904662965015
$ws.Range("B19").Value = 76.09584519098394
$ws.Range("B20").Value = 74.02422648832355
$ws.Range("B22").Value = 75.56528044338553
$ws.Range("B23").Value = 75.34079842925287
$ws.Range("B24").Value = 76.05801065791545
$ws.Range("B25").Value = 77.48682360591114
$ws.Range("B26").Value = 77.99155022377558
$ws.Range("B27").Value = 78.11627684164004
$ws.Range("B28").Value = 78.73139914350594
$ws.Range("B29").Value = 79.76701065283541
$ws.Range("B30").Value = 79.59669050896881
$ws.Range("B31").Value = 79.20092791936968
$ws.Range("B32").Value = 79.32452504203729
$ws.Range("B34").Value = 81.64151065029539
$ws.Range("B35").Value = 80.99215093802864
$ws.Range("B36").Value = 82.07303582949368
$ws.Range("B38").Value = 83.02601064775537
$ws.Range("B39").Value = 82.76081280575467
$ws.Range("B40").Value = 83.58546388348428
$ws.Range("B41").Value = 83.78290273255135
$ws.Range("B42").Value = 84.70339194001494
$ws.Range("B43").Value = 84.01218690468328
$ws.Range("B44").Value = 17.51876250323853
$ws.Range("B45").Value = 17.15139919430629
$ws.Range("B46").Value = 17.00612581217075
$ws.Range("B47").Value = 17.45925171070212
$ws.Range("B49").Value = 17.38993660115113
$ws.Range("B50").Value = 18.02082149261616
$ws.Range("B51").Value = 18.46226214001595
$ws.Range("B52").Value = 17.09412077277507
$ws.Range("B53").Value = 17.43823372229475
$ws.Range("B54").Value = 17.29599271522842
$ws.Range("B55").Value = 17.34059919024227
$ws.Range("B57").Value = 17.94501897393432
$ws.Range("B58").Value = 17.52335602416065
$ws.Range("B59").Value = 17.77164415364061
$ws.Range("B60").Value = 18.67202652794782
$ws.Range("B61").Value = 18.64956321850758
$ws.Range("B62").Value = 18.01945206732063
$ws.Range("B63").Value = 18.11655746282683
$ws.Range("B65").Value = 18.62219559154479
$ws.Range("B66").Value = 18.93080386489136
$ws.Range("B67").Value = 18.64584163495877
$ws.Range("B68").Value = 20.14
$ws.Range("B69").Value = 19.18000000000001
$ws.Range("B70").Value = 19.22000000000001
$ws.Range("B71").Value = 19.01000000000001
$ws.Range("B72").Value = 18.55
$ws.Range("B73").Value = 17.95
$ws.Range("B74").Value = 18.27
$ws.Range("B75").Value = 18.22
$ws.Range("B76").Value = 18.87
$ws.Range("B77").Value = 18.3
$ws.Range("B78").Value = 18.58
$ws.Range("B79").Value = 18.19000000000001
$ws.Range("B80").Value = 19.33
$ws.Range("B81").Value = 19.09
$ws.Range("B82").Value = 19.22
$ws.Range("B84").Value = 18.31999999999999
$ws.Range("B85").Value = 18.85000000000001
$ws.Range("B86").Value = 19.21000000000001
$ws.Range("B88").Value = 18.3
$ws.Range("B95").Value = 16.81999999999999
$ws.Range("B96").Value = 17.17999999999999
$ws.Range("B98").Value = 17.94
$ws.Range("B99").Value = 17.57000000000001
$ws.Range("B100").Value = 18.04000000000001
$ws.Range("B101").Value = 17.97
$ws.Range("B102").Value = 17.84999999999999
$ws.Range("B103").Value = 18.14
$ws.Range("B104").Value = 18.07000000000001
$ws.Range("B105").Value = 10.19
$ws.Range("B106").Value = 9.909999999999997
$ws.Range("B107").Value = 9.930000000000007
$ws.Range("B108").Value = 10.43000000000001
$ws.Range("B109").Value = 10.76000000000001
$ws.Range("B110").Value = 10.58
$ws.Range("B111").Value = 10.34
$ws.Range("B112").Value = 10.2
$ws.Range("B113").Value = 10.91000000000001
$ws.Range("B114").Value = 10.77
$ws.Range("B115").Value = 10.79000000000001
$ws.Range("B116").Value = 10.19
$ws.Range("B117").Value = 10.84
$ws.Range("B118").Value = 11.17999999999999
$ws.Range("B119").Value = 10.74000000000001
$ws.Range("B120").Value = 11.64999999999999
$ws.Range("B121").Value = 11.8
$ws.Range("B122").Value = 12.25
$ws.Range("B123").Value = 11.81999999999999
$ws.Range("B124").Value = 11.63
$ws.Range("B125").Value = 2.640000000000001
$ws.Range("B126").Value = 2.239999999999995
$ws.Range("B127").Value = 2.75
$ws.Range("B128").Value = 2.420000000000002
$ws.Range("B129").Value = 1.440000000000012
$ws.Range("B130").Value = 1.390000000000001
$ws.Range("B131").Value = 0.789999999999992
$ws.Range("B132").Value = -0.06000000000000227
$ws.Range("B133").Value = -0.1399999999999864
$ws.Range("B134").Value = 2.079999999999998
$ws.Range("B135").Value = 1.480000000000004
$ws.Range("B136").Value = 1.019999999999996
$ws.Range("B137").Value = 1.909999999999997
$ws.Range("B138").Value = 1.989999999999995
$ws.Range("B139").Value = 1.969999999999999
$ws.Range("B140").Value = 2
$ws.Range("B141").Value = 2.829999999999998
$ws.Range("B142").Value = 2.920000000000002
$ws.Range("B143").Value = 3.02000000000001
$ws.Range("B144").Value = 3.359999999999999
$ws.Range("B145").Value = 0.5
$ws.Range("B146").Value = 0.3599999999999994
$ws.Range("B147").Value = -0.03999999999999204
$ws.Range("B148").Value = 0.06999999999999318

# Add new row 149 (date 45748) with default revision value 0
$ws.Range("A149").Value = 45748
$ws.Range("B149").Value = 0

# Copy the date-cell formatting (style index used across column A) onto the new row
$ws.Range("A148").Copy()
$ws.Range("A149").PasteSpecial(-4122)
$excel.CutCopyMode = 0
